$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Junio de 2020 a las 19:11"

# Update country stats (refreshed data pull) and fix two pairs of swapped rows

# Row 4
$ws.Range("B4").Value = 2617847
$ws.Range("C4").Value = 21310
$ws.Range("D4").Value = 1082212
$ws.Range("E4").Value = 1407392
$ws.Range("G4").Value = 91
$ws.Range("H4").Value = 128243

# Row 7
$ws.Range("B7").Value = 548857
$ws.Range("C7").Value = 19280
$ws.Range("D7").Value = 321766
$ws.Range("E7").Value = 210606
$ws.Range("G7").Value = 382
$ws.Range("H7").Value = 16485

# Row 11
$ws.Range("B11").Value = 271982
$ws.Range("C11").Value = 4216
$ws.Range("D11").Value = 232210
$ws.Range("E11").Value = 34263
$ws.Range("G11").Value = 162
$ws.Range("H11").Value = 5509

# Row 12
$ws.Range("B12").Value = 240310
$ws.Range("C12").Value = 174
$ws.Range("D12").Value = 188891
$ws.Range("E12").Value = 16681
$ws.Range("G12").Value = 22
$ws.Range("H12").Value = 34738

# Row 22
$ws.Range("B22").Value = 103210
$ws.Range("C22").Value = 178
$ws.Range("D22").Value = 66152
$ws.Range("E22").Value = 28536
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 8522

# Row 39
$ws.Range("B39").Value = 41646
$ws.Range("C39").Value = 457
$ws.Range("D39").Value = 27066
$ws.Range("E39").Value = 13016
$ws.Range("G39").Value = 3
$ws.Range("H39").Value = 1564

# Row 49
$ws.Range("B49").Value = 25439
$ws.Range("C49").Value = 2
$ws.Range("E49").Value = 340
$ws.Range("G49").Value = 1
$ws.Range("H49").Value = 1735

# Row 50
$ws.Range("E50").Value = 5404
$ws.Range("G50").Value = 4
$ws.Range("H50").Value = 82

# Row 54
$ws.Range("E54").Value = 7778
$ws.Range("G54").Value = 12
$ws.Range("H54").Value = 178

# Row 63
$ws.Range("B63").Value = 13273
$ws.Range("C63").Value = 305
$ws.Range("D63").Value = 9371
$ws.Range("E63").Value = 3005
$ws.Range("G63").Value = 5
$ws.Range("H63").Value = 897

# Row 72
$ws.Range("B72").Value = 8853
$ws.Range("C72").Value = 7
$ws.Range("E72").Value = 466

# Row 80
$ws.Range("D80").Value = 1971
$ws.Range("E80").Value = 3956
$ws.Range("G80").Value = 2
$ws.Range("H80").Value = 143

# Row 105
$ws.Range("B105").Value = 2324
$ws.Range("C105").Value = 19
$ws.Range("D105").Value = 1911
$ws.Range("E105").Value = 405

# Row 128
$ws.Range("A128").Value = "Yemen"
$ws.Range("B128").Value = 1118
$ws.Range("C128").Value = 15
$ws.Range("D128").Value = 430
$ws.Range("E128").Value = 392
$ws.Range("H128").Value = 296

# Row 129
$ws.Range("A129").Value = "Letonia"
$ws.Range("B129").Value = 1116
$ws.Range("C129").Value = 1
$ws.Range("D129").Value = 932
$ws.Range("E129").Value = 154
$ws.Range("H129").Value = 30

# Row 143
$ws.Range("A143").Value = "Liberia"
$ws.Range("B143").Value = 768
$ws.Range("C143").Value = 39
$ws.Range("D143").Value = 298
$ws.Range("E143").Value = 436
$ws.Range("H143").Value = 34

# Row 144
$ws.Range("A144").Value = "Suazilandia"
$ws.Range("B144").Value = 745
$ws.Range("D144").Value = 370
$ws.Range("E144").Value = 367
$ws.Range("H144").Value = 8

# Row 205
$ws.Range("A205").Value = "Dominica"

# Row 206
$ws.Range("A206").Value = "Fiyi"

# Row 209
$ws.Range("A209").Value = "Islas Malvinas"

# Row 210
$ws.Range("A210").Value = "Groenlandia"
